# "first commit of 2023" — accept all tracked changes (insertions/deletions
# and the paragraph-mark deletion) left behind by N.D. Barber's 2021-11-11
# editing session, turning the Session->Lecture / 2021 header rewrite and
# its accompanying body-text edits into plain, un-tracked text.

$d = $word.ActiveDocument

# Make sure future edits (if any) aren't recorded as new tracked changes
# while we clean this revision history up.
$d.TrackRevisions = $false

# Accept every insertion/deletion (and paragraph-mark delete) still
# pending in the document - this is exactly what "Review > Accept All
# Changes" does in the Word UI.
$d.Revisions.AcceptAll()

Write-Output ("RemainingRevisions=" + $d.Revisions.Count)
